$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency ranking data (prices, volume %, and reordered rows).
# Each entry is: row, column, new value. Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
$updates = @(
    @(2, 4, '66.445.48'),
    @(2, 5, '  -0.72%  '),
    @(3, 4, '3.078.85'),
    @(3, 5, '  -1.20%  '),
    @(4, 5, '  +0.05%  '),
    @(5, 4, '573.97'),
    @(5, 5, '  -1.06%  '),
    @(6, 4, '170.42'),
    @(6, 5, '  -1.68%  '),
    @(7, 5, '  +0.13%  '),
    @(8, 4, '3.074.89'),
    @(8, 5, '  -1.20%  '),
    @(9, 4, '0.511'),
    @(9, 5, '  -1.99%  '),
    @(10, 4, '6.29'),
    @(10, 5, '  -2.08%  '),
    @(11, 4, '0.150'),
    @(11, 5, '  -2.32%  '),
    @(12, 5, '  -2.57%  '),
    @(13, 4, '0.0000238'),
    @(13, 5, '  -3.55%  '),
    @(14, 4, '35.66'),
    @(14, 5, '  -4.48%  '),
    @(15, 5, '  -0.85%  '),
    @(16, 4, '3.586.60'),
    @(16, 5, '  -1.23%  '),
    @(17, 4, '66.437.15'),
    @(17, 5, '  -0.63%  '),
    @(18, 4, '6.95'),
    @(18, 5, '  -2.67%  '),
    @(19, 4, '16.72'),
    @(19, 5, '  +2.14%  '),
    @(20, 4, '3.081.83'),
    @(20, 5, '  -1.06%  '),
    @(21, 4, '485.63'),
    @(21, 5, '  +1.67%  '),
    @(22, 2, 'Uniswap'),
    @(22, 3, 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'),
    @(22, 4, '7.81'),
    @(22, 5, '  -0.09%  '),
    @(23, 2, 'Polygon'),
    @(23, 3, 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @(23, 4, '0.684'),
    @(23, 5, '  -3.30%  '),
    @(24, 4, '82.71'),
    @(24, 5, '  -1.56%  '),
    @(25, 4, '12.67'),
    @(25, 5, '  -3.95%  '),
    @(26, 5, '  -3.83%  '),
    @(27, 2, 'Dai'),
    @(27, 3, 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @(27, 4, '1.00'),
    @(27, 5, '  +0.02%  '),
    @(28, 2, 'RenderToken'),
    @(28, 3, 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @(28, 4, '9.99'),
    @(28, 5, '  -3.71%  '),
    @(29, 4, '7.83'),
    @(29, 5, '  -1.57%  '),
    @(30, 5, '  -5.13%  '),
    @(31, 4, '2.58'),
    @(31, 5, '  -3.51%  '),
    @(32, 4, '27.79'),
    @(32, 5, '  -2.62%  '),
    @(33, 5, '  -3.60%  '),
    @(34, 4, '0.0₃0917'),
    @(34, 5, '  -3.36%  '),
    @(35, 4, '1.00'),
    @(35, 5, '  +0.08%  '),
    @(36, 4, '48.19'),
    @(36, 5, '  +2.40%  '),
    @(37, 2, 'Filecoin'),
    @(37, 3, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @(37, 4, '5.57'),
    @(37, 5, '  -4.65%  '),
    @(38, 2, 'Mantle'),
    @(38, 3, 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @(38, 4, '0.941'),
    @(38, 5, '  -3.26%  '),
    @(39, 2, 'OKB'),
    @(39, 3, 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
    @(39, 4, '48.93'),
    @(39, 5, '  -2.34%  '),
    @(40, 2, 'Kaspa'),
    @(40, 3, 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'),
    @(40, 4, '0.122'),
    @(40, 5, '  -1.48%  '),
    @(41, 2, 'TheGraph'),
    @(41, 3, 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'),
    @(41, 4, '0.303'),
    @(41, 5, '  -3.22%  '),
    @(42, 2, 'Stacks'),
    @(42, 3, 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'),
    @(42, 4, '1.95'),
    @(42, 5, '  -4.62%  '),
    @(43, 2, 'Cosmos'),
    @(43, 3, 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @(43, 4, '8.21'),
    @(43, 5, '  -3.98%  '),
    @(44, 2, 'Maker'),
    @(44, 3, 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @(44, 4, '2.770.89'),
    @(44, 5, '  -1.59%  '),
    @(45, 2, 'dogwifhat'),
    @(45, 3, 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'),
    @(45, 4, '2.55'),
    @(45, 5, '  -0.97%  '),
    @(46, 4, '366.62'),
    @(46, 5, '  -4.87%  '),
    @(47, 2, 'VeChain'),
    @(47, 3, 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @(47, 4, '0.0342'),
    @(47, 5, '  -3.04%  '),
    @(48, 2, 'Monero'),
    @(48, 3, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @(48, 4, '134.63'),
    @(48, 5, '  -0.97%  '),
    @(49, 2, 'USDe'),
    @(49, 3, 'https://coinranking.com/coin/exbfr2U-0+usde-usde'),
    @(49, 4, '1.00'),
    @(49, 5, '  -0.01%  '),
    @(50, 2, 'InjectiveProtocol'),
    @(50, 3, 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @(50, 4, '24.22'),
    @(50, 5, '  -3.43%  '),
    @(51, 2, 'ThetaToken'),
    @(51, 3, 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'),
    @(51, 4, '2.16'),
    @(51, 5, '  -2.05%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    # Force text format so numeric-looking strings (e.g. "1.00", "0.150") are
    # kept verbatim as text instead of being coerced into numeric cell values,
    # then restore the default (unstyled) cell style to match the original look.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
